$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" note under the title,
# keeping the (now blank) row itself.
$ws.Range("A2").Clear()

# The table used to report area for three census years (1989 / 2002 / 2014).
# Keep only the most recent (2014) column and drop the 1989/2002 columns
# (originally B:C); what was column D (2014 data) slides left into B.
$ws.Range("B1:C1").EntireColumn.Delete()

# Remove the now-unused spacer row that sat between the note and the
# "(sq. km)" sub-header.
$ws.Rows("3:3").Delete()

# Match the refreshed layout's taller rows, and extend formatting down to
# row 8 (matching the new used range of the sheet).
$ws.Range("A1:B8").RowHeight = 20.1
